# Updates the cryptos worksheet Price (D) and Volume(1h) (E) columns
# with the latest coinranking.com snapshot from this scheduled run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.895.66'
$ws.Range("E2").Value = '  -2.17%  '
$ws.Range("D3").Value = '1.834.66'
$ws.Range("E3").Value = '  -1.68%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.006'
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.46'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.005'
$ws.Range("E6").Value = '  +0.12%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4611'
$ws.Range("E7").Value = '  -1.50%  '
$ws.Range("E8").Value = '  -1.70%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07177'
$ws.Range("E9").Value = '  -2.78%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8799'
$ws.Range("E10").Value = '  -1.06%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07844'
$ws.Range("E11").Value = '  -1.17%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '19.63'
$ws.Range("E12").Value = '  -1.95%  '
$ws.Range("D13").Value = '1.857.88'
$ws.Range("E13").Value = '  -1.01%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.344'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.401'
$ws.Range("E15").Value = '  -2.97%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '88.39'
$ws.Range("E16").Value = '  -4.68%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.006'
$ws.Range("E17").Value = '  +0.09%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008759'
$ws.Range("E18").Value = '  -1.89%  '
$ws.Range("E19").Value = '  +0.06%  '
$ws.Range("D20").Value = '26.919.03'
$ws.Range("E20").Value = '  -2.20%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.49'
$ws.Range("E21").Value = '  -2.95%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.015'
$ws.Range("E22").Value = '  -2.82%  '
$ws.Range("E23").Value = '  -1.14%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.980'
$ws.Range("E24").Value = '  +5.30%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '150.86'
$ws.Range("E25").Value = '  -1.66%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '18.23'
$ws.Range("E26").Value = '  -1.63%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.991'
$ws.Range("E27").Value = '  -4.73%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '113.82'
$ws.Range("E28").Value = '  -2.72%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.960'
$ws.Range("E29").Value = '  -4.13%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.08842'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.127'
$ws.Range("E31").Value = '  +3.30%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7617'
$ws.Range("E32").Value = '  +0.76%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.467'
$ws.Range("E33").Value = '  -0.55%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.131'
$ws.Range("E34").Value = '  -2.92%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.657'
$ws.Range("E35").Value = '  +0.90%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.089'
$ws.Range("E36").Value = '  +0.72%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01936'
$ws.Range("E37").Value = '  -1.56%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.936'
$ws.Range("E38").Value = '  -1.76%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05141'
$ws.Range("E39").Value = '  -2.57%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.938'
$ws.Range("E40").Value = '  -3.39%  '
$ws.Range("E41").Value = '  -4.58%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1599'
$ws.Range("E42").Value = '  -2.87%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.312'
$ws.Range("E43").Value = '  -0.78%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4700'
$ws.Range("E44").Value = '  -3.66%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.005'
$ws.Range("E45").Value = '  +0.09%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.18'
$ws.Range("E46").Value = '  -1.76%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '102.84'
$ws.Range("E47").Value = '  -0.86%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.614'
$ws.Range("E48").Value = '  -2.68%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06091'
$ws.Range("E49").Value = '  -2.71%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '64.91'
$ws.Range("E50").Value = '  -1.41%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '36.39'
$ws.Range("E51").Value = '  -2.08%  '
